# Updates cryptos list prices/volumes per the Jun 24 2023 GitHub Actions data refresh.
# D-column (Price) values are apostrophe-prefixed so Excel stores them as
# text (quote-prefixed), matching the source inlineStr cells, instead of
# auto-converting the numeric-looking string to a Number -- which would
# drop significant trailing zeros (e.g. "0.9960" -> 0.996) or collapse the
# "." thousands separators (e.g. "30.572.63" -> 30572.63).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.572.63"
$ws.Range("E2").Value = "  -1.09%  "
$ws.Range("D3").Value = "'1.875.55"
$ws.Range("E3").Value = "  -1.56%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").Value = "'235.31"
$ws.Range("E5").Value = "  -4.63%  "
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("D7").Value = "'0.4865"
$ws.Range("E7").Value = "  -2.67%  "
$ws.Range("D8").Value = "'0.2895"
$ws.Range("E8").Value = "  -3.41%  "
$ws.Range("D9").Value = "'0.06641"
$ws.Range("E9").Value = "  -3.16%  "
$ws.Range("D10").Value = "'1.880.09"
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("D11").Value = "'16.53"
$ws.Range("E11").Value = "  -4.77%  "
$ws.Range("D12").Value = "'0.07227"
$ws.Range("E12").Value = "  -1.68%  "
$ws.Range("D13").Value = "'88.76"
$ws.Range("E13").Value = "  -3.56%  "
$ws.Range("D14").Value = "'4.983"
$ws.Range("E14").Value = "  -2.85%  "
$ws.Range("D15").Value = "'0.6494"
$ws.Range("E15").Value = "  -5.12%  "
$ws.Range("D16").Value = "'30.514.05"
$ws.Range("E16").Value = "  -1.20%  "
$ws.Range("D17").Value = "'0.000007817"
$ws.Range("E17").Value = "  -3.42%  "
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("D19").Value = "'12.92"
$ws.Range("E19").Value = "  -4.21%  "
$ws.Range("D20").Value = "'2.124.76"
$ws.Range("E20").Value = "  -1.25%  "
$ws.Range("D21").Value = "'1.004"
$ws.Range("E21").Value = "  +0.33%  "
$ws.Range("D22").Value = "'4.706"
$ws.Range("E22").Value = "  -3.67%  "
$ws.Range("D23").Value = "'190.72"
$ws.Range("E23").Value = "  +4.63%  "
$ws.Range("D24").Value = "'6.041"
$ws.Range("E24").Value = "  -1.16%  "
$ws.Range("D25").Value = "'9.251"
$ws.Range("E25").Value = "  -1.52%  "
$ws.Range("D26").Value = "'158.18"
$ws.Range("E26").Value = "  +2.44%  "
$ws.Range("D27").Value = "'18.22"
$ws.Range("E27").Value = "  -3.06%  "
$ws.Range("D28").Value = "'1.819"
$ws.Range("E28").Value = "  -7.10%  "
$ws.Range("D29").Value = "'1.406"
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("D30").Value = "'4.222"
$ws.Range("E30").Value = "  -3.97%  "
$ws.Range("D31").Value = "'0.08979"
$ws.Range("E31").Value = "  -0.16%  "
$ws.Range("D32").Value = "'3.905"
$ws.Range("E32").Value = "  -4.34%  "
$ws.Range("D33").Value = "'0.05122"
$ws.Range("E33").Value = "  -3.88%  "
$ws.Range("D34").Value = "'0.7196"
$ws.Range("E34").Value = "  -4.50%  "
$ws.Range("D35").Value = "'1.072"
$ws.Range("E35").Value = "  -6.40%  "
$ws.Range("E36").Value = "  -0.04%  "
$ws.Range("D37").Value = "'0.01801"
$ws.Range("E37").Value = "  -6.53%  "
$ws.Range("D38").Value = "'2.653"
$ws.Range("E38").Value = "  -2.97%  "
$ws.Range("D39").Value = "'0.9159"
$ws.Range("E39").Value = "  -2.68%  "
$ws.Range("D40").Value = "'2.032"
$ws.Range("E40").Value = "  -7.60%  "
$ws.Range("D41").Value = "'0.4349"
$ws.Range("E41").Value = "  -1.36%  "
$ws.Range("D42").Value = "'104.27"
$ws.Range("E42").Value = "  -2.07%  "
$ws.Range("D43").Value = "'0.9960"
$ws.Range("E43").Value = "  -0.44%  "
$ws.Range("D44").Value = "'5.693"
$ws.Range("E44").Value = "  -3.18%  "
$ws.Range("D45").Value = "'0.1325"
$ws.Range("E45").Value = "  -3.04%  "
$ws.Range("D46").Value = "'7.294"
$ws.Range("E46").Value = "  -6.55%  "
$ws.Range("D47").Value = "'0.4021"
$ws.Range("E47").Value = "  +2.08%  "
$ws.Range("D48").Value = "'0.05822"
$ws.Range("E48").Value = "  -0.51%  "
$ws.Range("D49").Value = "'8.619"
$ws.Range("E49").Value = "  -0.14%  "
$ws.Range("D50").Value = "'1.397"
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("D51").Value = "'33.09"
$ws.Range("E51").Value = "  -1.62%  "
